$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (Strike#) column values for rows 2-17 (column G)
$kValues = @{
    2  = 6
    3  = 2
    4  = 0
    5  = 1
    6  = 4
    7  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 4
    13 = 3
    14 = 2
    15 = 2
    16 = 2
    17 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
